$d = $word.ActiveDocument

# Insert a new first paragraph containing "Testing" followed by a
# placeholder character and a paragraph break. The placeholder lets us
# anchor the "_GoBack" bookmark at a true interior position (so it ends
# up fully contained in the new paragraph, after the "Testing" run,
# rather than snapping across the paragraph-mark boundary).
$r = $d.Range(0, 0)
$r.InsertBefore("TestingX`r")

# Re-create the "_GoBack" bookmark around the placeholder character.
# Adding a bookmark with the same name as the existing one removes the
# old one (bookmark names are unique), which takes care of deleting the
# bookmark that used to sit at the end of the picture paragraph.
$bmRange = $d.Range(7, 8)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character through the bookmark's own range so
# the now-empty bookmark collapses in place (right after "Testing",
# still inside the new first paragraph) instead of being deleted.
$bm = $d.Bookmarks("_GoBack")
$delRange = $bm.Range
$delRange.Text = ""
